$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the order date (A1) by one month: 2024-04-24 -> 2024-05-24
$ws.Range("A1").Value = 45436

# Update prices for the two "PATAS PLASTICAS" rows
$ws.Range("D29").Value = 1230
$ws.Range("D30").Value = 1290
